$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3359.4
$ws.Range("I34").Value = 3359.4
$ws.Range("K34").Value = 3359.4
$ws.Range("M34").Value = -3156.4
$ws.Range("H36").Value = 3359.4
$ws.Range("I36").Value = 3359.4
$ws.Range("K36").Value = 3359.4
$ws.Range("M36").Value = -2644.4
$ws.Range("H107").Value = 1845.8
$ws.Range("I107").Value = 1290
$ws.Range("J107").Value = 2401.6
$ws.Range("K107").Value = 1290
$ws.Range("L107").Value = 2401.6
$ws.Range("M107").Value = 630
$ws.Range("N107").Value = -6241.6
$ws.Range("H111").Value = 722.2857
$ws.Range("J111").Value = 950
$ws.Range("L111").Value = 2850
$ws.Range("N111").Value = -8984

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1533.25
$ws.Range("I5").Value = 66.5
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 66.5
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 45.5
$ws.Range("N5").Value = -3224
$ws.Range("H35").Value = 1524.2
$ws.Range("I35").Value = 1524.2
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1524.2
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1118.2
$ws.Range("N35").ClearContents()
$ws.Range("H50").Value = 15072
$ws.Range("I50").Value = 5147.75
$ws.Range("J50").Value = 24996.25
$ws.Range("K50").Value = 5147.75
$ws.Range("L50").Value = 24996.25
$ws.Range("M50").Value = -4433.75
$ws.Range("N50").Value = -26424.25
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H60").Value = 31017
$ws.Range("I60").Value = 2051
$ws.Range("K60").Value = 2051
$ws.Range("M60").Value = -1318
$ws.Range("H122").Value = 4323.375
$ws.Range("I122").Value = 4323.375
$ws.Range("K122").Value = 12970.125
$ws.Range("M122").Value = -10520.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1533.25
$ws.Range("I4").Value = 66.5
$ws.Range("J4").Value = 3000
$ws.Range("K4").Value = 66.5
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 48.5
$ws.Range("N4").Value = -3230
$ws.Range("H36").Value = 6348
$ws.Range("I36").Value = 6666.3335
$ws.Range("K36").Value = 6666.3335
$ws.Range("M36").Value = -6132.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 354
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -1200
$ws.Range("H31").Value = 8767.049999999999
$ws.Range("I31").Value = 2830.2
$ws.Range("J31").Value = 26577.6
$ws.Range("K31").Value = 2830.2
$ws.Range("L31").Value = 26577.6
$ws.Range("M31").Value = -2535.2
$ws.Range("N31").Value = -27167.6
$ws.Range("H34").Value = 8767.049999999999
$ws.Range("I34").Value = 2830.2
$ws.Range("J34").Value = 26577.6
$ws.Range("K34").Value = 2830.2
$ws.Range("L34").Value = 26577.6
$ws.Range("M34").Value = -2628.2
$ws.Range("N34").Value = -26981.6
$ws.Range("H35").Value = 3578.5
$ws.Range("I35").Value = 868
$ws.Range("J35").Value = 8999.5
$ws.Range("K35").Value = 868
$ws.Range("L35").Value = 8999.5
$ws.Range("M35").Value = -574
$ws.Range("N35").Value = -9587.5
$ws.Range("H36").Value = 4696.4
$ws.Range("I36").Value = 4696.4
$ws.Range("K36").Value = 4696.4
$ws.Range("M36").Value = -4308.4
$ws.Range("H40").Value = 4696.4
$ws.Range("I40").Value = 4696.4
$ws.Range("K40").Value = 4696.4
$ws.Range("M40").Value = -4536.4
$ws.Range("H42").Value = 19598.8
$ws.Range("I42").Value = 8997
$ws.Range("J42").Value = 22249.25
$ws.Range("K42").Value = 8997
$ws.Range("L42").Value = 22249.25
$ws.Range("M42").Value = -8404
$ws.Range("N42").Value = -23435.25
$ws.Range("H44").Value = 29997.8
$ws.Range("I44").Value = 29996.334
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 29996.334
$ws.Range("L44").Value = 30000
$ws.Range("M44").Value = -29554.334
$ws.Range("N44").Value = -30884
$ws.Range("H55").Value = 8666.333000000001
$ws.Range("I55").Value = 8666.333000000001
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 8666.333000000001
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -8351.333000000001
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H114").Value = 769.1667
$ws.Range("J114").Value = 467.5
$ws.Range("L114").Value = 1402.5
$ws.Range("N114").Value = -7910.5
$ws.Range("H129").Value = 2294.6
$ws.Range("I129").Value = 368.25
$ws.Range("K129").Value = 1104.75
$ws.Range("M129").Value = 3895.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 8896.571
$ws.Range("I41").Value = 9212.666999999999
$ws.Range("K41").Value = 9212.666999999999
$ws.Range("M41").Value = -8857.666999999999
$ws.Range("H126").Value = 11999.25
$ws.Range("I126").Value = 11999.25
$ws.Range("K126").Value = 35997.75
$ws.Range("M126").Value = -33527.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 2083.3333
$ws.Range("I32").Value = 2083.3333
$ws.Range("K32").Value = 2083.3333
$ws.Range("M32").Value = -1766.3333
$ws.Range("H43").Value = 155634.1
$ws.Range("J43").Value = 155634.1
$ws.Range("L43").Value = 155634.1
$ws.Range("N43").Value = -156020.1
$ws.Range("H54").Value = 19082.334
$ws.Range("J54").Value = 19082.334
$ws.Range("L54").Value = 19082.334
$ws.Range("N54").Value = -20370.334
$ws.Range("H57").Value = 20000
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 20000
$ws.Range("N57").Value = -21132
$ws.Range("H58").Value = 31736.666
$ws.Range("I58").Value = 3210
$ws.Range("K58").Value = 3210
$ws.Range("M58").Value = -2950
$ws.Range("H68").Value = 2999
$ws.Range("I68").Value = 2999
$ws.Range("K68").Value = 2999
$ws.Range("M68").Value = -2250
$ws.Range("H71").Value = 2999
$ws.Range("I71").Value = 2999
$ws.Range("K71").Value = 14995
$ws.Range("M71").Value = -11251
$ws.Range("H93").Value = 965.6667
$ws.Range("I93").Value = 998
$ws.Range("J93").Value = 949.5
$ws.Range("K93").Value = 998
$ws.Range("L93").Value = 949.5
$ws.Range("M93").Value = 250
$ws.Range("N93").Value = -3445.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 37000000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H51").Value = 27657.857
$ws.Range("J51").Value = 25629.8
$ws.Range("L51").Value = 25629.8
$ws.Range("N51").Value = -26649.8
$ws.Range("H52").Value = 24772.5
$ws.Range("J52").Value = 24772.5
$ws.Range("L52").Value = 24772.5
$ws.Range("N52").Value = -25224.5
$ws.Range("H62").Value = 926
$ws.Range("I62").Value = 926
$ws.Range("K62").Value = 926
$ws.Range("M62").Value = -302
$ws.Range("H65").Value = 926
$ws.Range("I65").Value = 926
$ws.Range("K65").Value = 4630
$ws.Range("M65").Value = -1510
$ws.Range("H107").Value = 426
$ws.Range("I107").Value = 426
$ws.Range("K107").Value = 1278
$ws.Range("M107").Value = 642
